# Apply data-refresh edits (generated output update) to the 广州-漫展信息 workbook.
# The commit corresponds to a re-generated data pull: a handful of "views/likes"
# counters (column F) were bumped, and one venue address (column D) was renamed
# from "广州国际医药港" to "广州健康方舟". The same source rows are mirrored on
# sheet "全部类型" (which aggregates 展览 + 演出 + 本地生活), so every change is
# applied there too.
#
# NOTE: named parameter binding (-Foo bar) does not reliably pass values through
# in this runtime, so the helper below is called positionally.

$wb = $excel.ActiveWorkbook

function Set-CellValue($SheetName, $CellRef, $Value) {
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($CellRef).Value = $Value
}

# ---- Sheet "展览" (exhibitions) ----
Set-CellValue "展览" "F2"  831
Set-CellValue "展览" "F3"  13652
Set-CellValue "展览" "F4"  13450
Set-CellValue "展览" "F7"  36
Set-CellValue "展览" "F12" 738
Set-CellValue "展览" "F13" 2125
Set-CellValue "展览" "D14" "东沙大道16号 广州健康方舟"
Set-CellValue "展览" "F14" 66
Set-CellValue "展览" "F20" 363
Set-CellValue "展览" "F22" 495
Set-CellValue "展览" "F24" 64

# ---- Sheet "演出" (performances) ----
Set-CellValue "演出" "F6" 154
Set-CellValue "演出" "F7" 1360
Set-CellValue "演出" "F9" 22

# ---- Sheet "本地生活" (local life) ----
Set-CellValue "本地生活" "F2" 215

# ---- Sheet "全部类型" (all types - combined view of the above) ----
Set-CellValue "全部类型" "F2"  215
Set-CellValue "全部类型" "F3"  831
Set-CellValue "全部类型" "F4"  13652
Set-CellValue "全部类型" "F5"  13450
Set-CellValue "全部类型" "F8"  36
Set-CellValue "全部类型" "F13" 738
Set-CellValue "全部类型" "F16" 2125
Set-CellValue "全部类型" "D17" "东沙大道16号 广州健康方舟"
Set-CellValue "全部类型" "F17" 66
Set-CellValue "全部类型" "F27" 363
Set-CellValue "全部类型" "F29" 495
Set-CellValue "全部类型" "F31" 154
Set-CellValue "全部类型" "F32" 1360
Set-CellValue "全部类型" "F34" 22
Set-CellValue "全部类型" "F35" 64
